# ------------------------------------------------------------------
# Apply the commit:
#   - insert a new first sheet "Player Info" with player bio data
#   - rename column D ("ODI Batting") / column B ("ODI Bowling") header
#     from MATCH_CARD_LINK -> MATCH_CODE, and replace the full scorecard
#     URL values with just the numeric match code
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Insert a new worksheet named "Player Info" before "ODI Batting".
# NOTE: worksheet handles obtained before a structural change such as
# Worksheets.Add() become stale (they track position, not identity), so
# every sheet reference used below is (re-)fetched *after* the insert.
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$newSheet = $wb.Worksheets.Add($battingSheetForInsert)
$newSheet.Name = "Player Info"

# Re-fetch stable references now that the sheet collection has changed.
$infoSheet = $wb.Worksheets.Item("Player Info")
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Match the page margins used on the other sheets (0.75/0.75/1/1/0.5/0.5 in)
$infoSheet.PageSetup.LeftMargin = 0.75 * 72
$infoSheet.PageSetup.RightMargin = 0.75 * 72
$infoSheet.PageSetup.TopMargin = 1 * 72
$infoSheet.PageSetup.BottomMargin = 1 * 72
$infoSheet.PageSetup.HeaderMargin = 0.5 * 72
$infoSheet.PageSetup.FooterMargin = 0.5 * 72

# Header row (bold / thin-bordered / center+top aligned, matching the
# header style already used on the "ODI Batting" / "ODI Bowling" sheets)
$infoSheet.Cells.Item(1, 1).Value = "ID"
$infoSheet.Cells.Item(1, 2).Value = "NAME"
$infoSheet.Cells.Item(1, 3).Value = "BATTING_HAND"
$infoSheet.Cells.Item(1, 4).Value = "BOWL_STYLE"

$headerRange = $infoSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row
$infoSheet.Cells.Item(2, 1).Value = "'4404"
$infoSheet.Cells.Item(2, 2).Value = "Deepak Jagbir Hooda"
$infoSheet.Cells.Item(2, 3).Value = "Right Handed"
$infoSheet.Cells.Item(2, 4).Value = "Right Arm Off Break"

$infoSheet.Range("A1").Select()

# ---- 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4533"
    3  = "4535"
    4  = "4621"
    5  = "4623"
    6  = "4624"
    7  = "4637"
    8  = "4640"
    9  = "4643"
    10 = "4673"
    11 = "4676"
}
foreach ($row in $battingCodes.Keys) {
    $battingSheet.Cells.Item($row, 4).Value = "'" + $battingCodes[$row]
}

# ---- 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4535"
    3 = "4621"
    4 = "4623"
    5 = "4624"
    6 = "4640"
    7 = "4643"
}
foreach ($row in $bowlingCodes.Keys) {
    $bowlingSheet.Cells.Item($row, 2).Value = "'" + $bowlingCodes[$row]
}
